# Change the table placeholder syntax from "${table:name.column}" to
# "${table:name:column}" so the processor can support multi-dimensional
# key names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Tables")

$ws.Range("B5").Value = '${table:ages:name}'
$ws.Range("C5").Value = '${table:ages:age}'
$ws.Range("B8").Value = '${table:hours:name}'
$ws.Range("C8").Value = '${table:hours:days}'
